$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "vfdvfd"
$ws.Range("B2").Value = 17
# C2 keeps its original text/inline-string type ("1000" is numeric-looking,
# so a leading apostrophe forces Excel to store it as text rather than a
# number); reset style afterwards so it doesn't pick up a quote-prefix style.
$ws.Range("C2").Value = "'1000"
$ws.Range("C2").Style = "Normal"

# Row 3
$ws.Range("A3").Value = "ekjnerk"
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = "'1200"
$ws.Range("C3").Style = "Normal"
